$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "66.553.68"
$ws.Range("E2").Value = "  +3.88%  "

$ws.Range("D3").Value = "3.485.02"
$ws.Range("E3").Value = "  +2.31%  "

$ws.Range("E4").Value = "  +0.04%  "

Set-TextValue $ws.Range("D5") "591.83"
$ws.Range("E5").Value = "  +3.62%  "

Set-TextValue $ws.Range("D6") "168.96"
$ws.Range("E6").Value = "  +3.98%  "

$ws.Range("D8").Value = "3.482.80"
$ws.Range("E8").Value = "  +2.23%  "

Set-TextValue $ws.Range("D9") "0.591"
$ws.Range("E9").Value = "  +8.06%  "

$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").Value = "  +6.88%  "

Set-TextValue $ws.Range("D12") "0.436"
$ws.Range("E12").Value = "  +4.16%  "

$ws.Range("D13").Value = "4.085.20"
$ws.Range("E13").Value = "  +2.41%  "

Set-TextValue $ws.Range("D14") "0.134"
$ws.Range("E14").Value = "  -0.32%  "

Set-TextValue $ws.Range("D15") "28.14"
$ws.Range("E15").Value = "  +4.75%  "

$ws.Range("E16").Value = "  +3.44%  "

$ws.Range("D17").Value = "66.596.18"
$ws.Range("E17").Value = "  +3.91%  "

$ws.Range("D18").Value = "3.480.49"
$ws.Range("E18").Value = "  +1.40%  "

$ws.Range("E19").Value = "  +3.36%  "

Set-TextValue $ws.Range("D20") "14.01"
$ws.Range("E20").Value = "  +4.01%  "

Set-TextValue $ws.Range("D21") "392.01"
$ws.Range("E21").Value = "  +5.43%  "

$ws.Range("E22").Value = "  +1.68%  "

Set-TextValue $ws.Range("D23") "72.97"
$ws.Range("E23").Value = "  +4.15%  "

Set-TextValue $ws.Range("D24") "0.999"

Set-TextValue $ws.Range("D25") "0.534"
$ws.Range("E25").Value = "  +4.79%  "

Set-TextValue $ws.Range("D26") "0.0000121"
$ws.Range("E26").Value = "  +6.00%  "

Set-TextValue $ws.Range("D27") "10.32"
$ws.Range("E27").Value = "  +8.74%  "

Set-TextValue $ws.Range("D28") "0.180"
$ws.Range("E28").Value = "  +1.44%  "

Set-TextValue $ws.Range("D29") "1.00"
$ws.Range("E29").Value = "  +0.09%  "

Set-TextValue $ws.Range("D30") "6.31"
$ws.Range("E30").Value = "  +4.09%  "

$ws.Range("E31").Value = "  +4.67%  "

Set-TextValue $ws.Range("D32") "2.06"
$ws.Range("E32").Value = "  +3.55%  "

Set-TextValue $ws.Range("D33") "23.58"
$ws.Range("E33").Value = "  +3.81%  "

$ws.Range("E34").Value = "  +5.86%  "

Set-TextValue $ws.Range("D35") "0.999"
$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +9.01%  "

Set-TextValue $ws.Range("D37") "161.34"
$ws.Range("E37").Value = "  +1.12%  "

Set-TextValue $ws.Range("D38") "0.893"
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("E39").Value = "  +6.83%  "

$ws.Range("E40").Value = "  +5.63%  "

Set-TextValue $ws.Range("D41") "0.0745"
$ws.Range("E41").Value = "  +3.43%  "

Set-TextValue $ws.Range("D42") "26.49"
$ws.Range("E42").Value = "  +2.68%  "

$ws.Range("E43").Value = "  +6.74%  "

Set-TextValue $ws.Range("D44") "26.86"
$ws.Range("E44").Value = "  +3.87%  "

Set-TextValue $ws.Range("D45") "43.15"
$ws.Range("E45").Value = "  +1.22%  "

$ws.Range("D46").Value = "2.765.12"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("E47").Value = "  +3.45%  "

Set-TextValue $ws.Range("D48") "2.48"
$ws.Range("E48").Value = "  +3.87%  "

Set-TextValue $ws.Range("D49") "345.52"
$ws.Range("E49").Value = "  +5.08%  "

$ws.Range("E50").Value = "  +5.08%  "

Set-TextValue $ws.Range("D51") "33.91"
$ws.Range("E51").Value = "  +12.81%  "
